# Go: improve the way formula are created (WIP)
#
# Adds a "Basic arithmetic" worked example below the existing type examples:
#   A8: label "Basic arithmetic" (cell already carries the bold header style)
#   A9: =B9+C9   (evaluates to 2)
#   B9: 1
#   C9: 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Basic arithmetic"

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("A9").Formula = "=B9+C9"

[void]$ws.Range("A9").Select()
